# The workbook's single sheet is dated "Opdateret d. 02-12-2025" and is
# referenced both as the sheet's tab name and inside the workbook-level
# defined name "Kiropraktorsystemer". The data was refreshed on
# 05-12-2025, so the sheet (and therefore the defined name that points at
# it) needs to be renamed to match the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldName = "Opdateret d. 02-12-2025"
$newName = "Opdateret d. 05-12-2025"

if ($ws.Name -eq $oldName) {
    $ws.Name = $newName
} else {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $newName
}
